$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 5; $row++) {
    $ws.Cells.Item($row, 2).Value = "Yes"   # Column B: Did Harvest Occur?
    $ws.Cells.Item($row, 6).Value = "Na"    # Column F: Species
    $ws.Cells.Item($row, 10).Value = 1      # Column J: Unknown Sex Count
}
